{"js": "const body = context.document.body;\n\n// Map of old cell text -> new cell text (three-digit \u00f7 one-digit practice answers).\nconst replacements = [\n  [\"421\u00f74=105, 1\", \"482\u00f79=53, 5\"],\n  [\"477\u00f78=59, 5\", \"738\u00f78=92, 2\"],\n  [\"936\u00f76=156, 0\", \"741\u00f72=370, 1\"],\n  [\"177\u00f73=59, 0\", \"931\u00f75=186, 1\"],\n  [\"908\u00f75=181, 3\", \"893\u00f73=297, 2\"],\n  [\"104\u00f74=26, 0\", \"877\u00f79=97, 4\"],\n  [\"110\u00f73=36, 2\", \"501\u00f72=250, 1\"],\n  [\"142\u00f75=28, 2\", \"317\u00f72=158, 1\"],\n  [\"653\u00f74=163, 1\", \"351\u00f79=39, 0\"],\n  [\"978\u00f78=122, 2\", \"419\u00f78=52, 3\"],\n  [\"754\u00f74=188, 2\", \"882\u00f77=126, 0\"],\n  [\"548\u00f75=109, 3\", \"846\u00f79=94, 0\"],\n  [\"318\u00f76=53, 0\", \"338\u00f72=169, 0\"],\n  [\"496\u00f76=82, 4\", \"895\u00f74=223, 3\"],\n  [\"450\u00f75=90, 0\", \"781\u00f78=97, 5\"],\n  [\"865\u00f76=144, 1\", \"257\u00f75=51, 2\"],\n  [\"653\u00f73=217, 2\", \"282\u00f74=70, 2\"],\n  [\"826\u00f76=137, 4\", \"645\u00f73=215, 0\"],\n  [\"541\u00f74=135, 1\", \"496\u00f78=62, 0\"],\n  [\"581\u00f72=290, 1\", \"785\u00f79=87, 2\"],\n  [\"330\u00f75=66, 0\", \"636\u00f77=90, 6\"],\n  [\"700\u00f76=116, 4\", \"771\u00f72=385, 1\"],\n  [\"650\u00f76=108, 2\", \"877\u00f72=438, 1\"],\n  [\"904\u00f79=100, 4\", \"912\u00f79=101, 3\"],\n  [\"502\u00f74=125, 2\", \"914\u00f72=457, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Map of old cell text -> new cell text (three-digit \u00f7 one-digit practice answers).\n$replacements = @(\n    @(\"421\u00f74=105, 1\", \"482\u00f79=53, 5\"),\n    @(\"477\u00f78=59, 5\", \"738\u00f78=92, 2\"),\n    @(\"936\u00f76=156, 0\", \"741\u00f72=370, 1\"),\n    @(\"177\u00f73=59, 0\", \"931\u00f75=186, 1\"),\n    @(\"908\u00f75=181, 3\", \"893\u00f73=297, 2\"),\n    @(\"104\u00f74=26, 0\", \"877\u00f79=97, 4\"),\n    @(\"110\u00f73=36, 2\", \"501\u00f72=250, 1\"),\n    @(\"142\u00f75=28, 2\", \"317\u00f72=158, 1\"),\n    @(\"653\u00f74=163, 1\", \"351\u00f79=39, 0\"),\n    @(\"978\u00f78=122, 2\", \"419\u00f78=52, 3\"),\n    @(\"754\u00f74=188, 2\", \"882\u00f77=126, 0\"),\n    @(\"548\u00f75=109, 3\", \"846\u00f79=94, 0\"),\n    @(\"318\u00f76=53, 0\", \"338\u00f72=169, 0\"),\n    @(\"496\u00f76=82, 4\", \"895\u00f74=223, 3\"),\n    @(\"450\u00f75=90, 0\", \"781\u00f78=97, 5\"),\n    @(\"865\u00f76=144, 1\", \"257\u00f75=51, 2\"),\n    @(\"653\u00f73=217, 2\", \"282\u00f74=70, 2\"),\n    @(\"826\u00f76=137, 4\", \"645\u00f73=215, 0\"),\n    @(\"541\u00f74=135, 1\", \"496\u00f78=62, 0\"),\n    @(\"581\u00f72=290, 1\", \"785\u00f79=87, 2\"),\n    @(\"330\u00f75=66, 0\", \"636\u00f77=90, 6\"),\n    @(\"700\u00f76=116, 4\", \"771\u00f72=385, 1\"),\n    @(\"650\u00f76=108, 2\", \"877\u00f72=438, 1\"),\n    @(\"904\u00f79=100, 4\", \"912\u00f79=101, 3\"),\n    @(\"502\u00f74=125, 2\", \"914\u00f72=457, 0\"),\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
